$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73-153 down to 74-154.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record.
$ws.Range("A73").Value = 10
$ws.Range("B73").Value = "Vega Modelo de Temuco"
$ws.Range("C73").Value = "La Araucanía"
$ws.Range("D73").Value = 44494
$ws.Range("E73").Value = 9
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100101
$ws.Range("H73").Value = "Berries"
$ws.Range("I73").Value = 100112025
$ws.Range("J73").Value = "Frutilla"
$ws.Range("K73").Value = "Sin especificar"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 550
$ws.Range("N73").Value = 10000
$ws.Range("O73").Value = 10000
$ws.Range("P73").Value = 10000
$ws.Range("Q73").Value = "$/bandeja 7 kilos"
$ws.Range("R73").Value = "Provincia de Melipilla"
$ws.Range("S73").Value = 1429
$ws.Range("T73").Value = 7
